$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above row 290 (shifts old rows 290-302 down to 294-306)
$ws.Rows.Item(290).Insert()
$ws.Rows.Item(290).Insert()
$ws.Rows.Item(290).Insert()
$ws.Rows.Item(290).Insert()

# New data rows to populate at 290-293
$newRows = @(
    @{Row=290; D=44615; K="August Red"; L="Primera"; M=22; N=310000; O=320000; P=315000; S=750},
    @{Row=291; D=44615; K="August Red"; L="Segunda"; M=18; N=270000; O=280000; P=275000; S=655},
    @{Row=292; D=44615; K="Venus";      L="Especial"; M=24; N=320000; O=330000; P=325000; S=774},
    @{Row=293; D=44615; K="Venus";      L="Primera";  M=18; N=300000; O=310000; P=305000; S=726}
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 8
    $ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100103
    $ws.Cells.Item($row, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($row, 9).Value = 100103006
    $ws.Cells.Item($row, 10).Value = "Nectarín"
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "$/bins (420 kilos)"
    $ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 420
}
